# "Dropdowns for global objects"
#
# The sheet originally held two lookup columns (A: Es-style module keys,
# B: Ep-style module keys) feeding two dropdown lists. This change adds a
# second, "global objects" pair of dropdown columns: the existing A/B
# columns are duplicated into new C/D columns (keeping their original
# headers + row data untouched), and the original A/B headers are
# replaced with the new "EP.*" global-object header keys. Both header
# rows are bolded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing A1:B11 lookup table into C1:D11 - this keeps the
# original "Functions.Es/EpOpenModule.module" columns (with all their
# header + row values) intact, just shifted over.
$ws.Range("A1:B11").Copy() | Out-Null
$ws.Range("C1").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

# Re-purpose the original A/B columns as the new "global objects" dropdown
# pair - only the header row text changes, row data (2-11) is shared.
$ws.Range("A1").Value = "EP.EsOpenModule.module"
$ws.Range("B1").Value = "EP.OpenModule.module"

# Bold both header rows (new A1:B1 and the carried-over C1:D1).
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("C1:D1").Font.Bold = $true

# Best-effort column widths to match the post-edit column layout.
$ws.Columns.Item(1).ColumnWidth = 30.8307
$ws.Columns.Item(2).ColumnWidth = 31.0534
$ws.Columns.Item(3).ColumnWidth = 30.8307
$ws.Columns.Item(4).ColumnWidth = 29.1667
$ws.Columns.Item(5).ColumnWidth = 29.1667

# Restore the cursor to the top of the new primary dropdown column.
$ws.Range("B2").Select() | Out-Null
